$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$comment = "From the AerChem side there is interest in Emon reffclws. This variable is identified as the already available PEXTRA variable with the table 126 grib code 126021, i.e. proposing to add reffclws as 21.126 to ifspar.json.  Note that this variable is not requested by CMIP6 AerChem, and that reffclws not occurs in any CMIP6 data request of the experiments in which any EC-Earth3* configuration participates. See further #564."

# Row 44: CF3hr / reffclws
$ws.Range("A44").Value = "CF3hr"
$ws.Range("B44").Value = "reffclws"
$ws.Range("H44").Value = $comment
$ws.Range("I44").Value = "Thomas"

# Row 45: Esubhr / reffclws
$ws.Range("A45").Value = "Esubhr"
$ws.Range("B45").Value = "reffclws"
$ws.Range("H45").Value = $comment
$ws.Range("I45").Value = "Thomas"

# Wrap text style for column B cells of the new rows (matches new cellXf with wrapText=true)
$ws.Range("B44:B45").WrapText = $true

$ws.Range("A44").Select()
